$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be misread as a number by Excel
# (e.g. "323.68", "0.9999", trailing-zero decimals like "48.30") are forced
# to Text format first so the literal string is preserved verbatim.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D12",
    "D14",
    "D15",
    "D16",
    "D19",
    "D20",
    "D22",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D46",
    "D47",
    "D48",
    "D49",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices, volume deltas, and the swapped
# PancakeSwap / NEARProtocol row 48<->49 content) per the diff.
$ws.Range("D2").Value = '27.364.77'
$ws.Range("E2").Value = '  -3.93%  '
$ws.Range("D3").Value = '1.861.97'
$ws.Range("E3").Value = '  -4.76%  '
$ws.Range("E4").Value = '  -1.01%  '
$ws.Range("D5").Value = '323.68'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").Value = '0.4531'
$ws.Range("E7").Value = '  -5.57%  '
$ws.Range("D8").Value = '0.3871'
$ws.Range("E8").Value = '  -5.15%  '
$ws.Range("D9").Value = '48.30'
$ws.Range("E9").Value = '  -10.51%  '
$ws.Range("D10").Value = '0.07910'
$ws.Range("E10").Value = '  -6.82%  '
$ws.Range("E11").Value = '  -3.39%  '
$ws.Range("D12").Value = '21.42'
$ws.Range("E12").Value = '  -4.58%  '
$ws.Range("D13").Value = '1.861.78'
$ws.Range("E13").Value = '  -5.24%  '
$ws.Range("D14").Value = '5.904'
$ws.Range("E14").Value = '  -4.24%  '
$ws.Range("D15").Value = '7.152'
$ws.Range("E15").Value = '  -5.64%  '
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("E17").Value = '  -3.80%  '
$ws.Range("E18").Value = '  -5.42%  '
$ws.Range("D19").Value = '0.06503'
$ws.Range("E19").Value = '  -1.89%  '
$ws.Range("D20").Value = '17.13'
$ws.Range("E20").Value = '  -7.36%  '
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").Value = '5.529'
$ws.Range("E22").Value = '  -5.31%  '
$ws.Range("D23").Value = '27.363.49'
$ws.Range("E23").Value = '  -4.01%  '
$ws.Range("D24").Value = '10.89'
$ws.Range("E24").Value = '  -4.91%  '
$ws.Range("D25").Value = '2.273'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = '2.078.89'
$ws.Range("E26").Value = '  -5.31%  '
$ws.Range("D27").Value = '152.94'
$ws.Range("E27").Value = '  -2.35%  '
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("D29").Value = '2.064'
$ws.Range("E29").Value = '  -5.13%  '
$ws.Range("D30").Value = '5.487'
$ws.Range("E30").Value = '  -5.72%  '
$ws.Range("D31").Value = '120.94'
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("E32").Value = '  +3.37%  '
$ws.Range("D33").Value = '0.09319'
$ws.Range("E33").Value = '  -3.69%  '
$ws.Range("D34").Value = '0.9351'
$ws.Range("E34").Value = '  -5.20%  '
$ws.Range("D35").Value = '3.604'
$ws.Range("E35").Value = '  -2.29%  '
$ws.Range("D36").Value = '5.276'
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("E37").Value = '  -4.08%  '
$ws.Range("D38").Value = '1.223'
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").Value = '0.06001'
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").Value = '8.254'
$ws.Range("E40").Value = '  -9.57%  '
$ws.Range("D41").Value = '0.9997'
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").Value = '0.5913'
$ws.Range("E42").Value = '  -5.08%  '
$ws.Range("D43").Value = '0.1892'
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("E44").Value = '  -9.31%  '
$ws.Range("E45").Value = '  -5.13%  '
$ws.Range("D46").Value = '0.5629'
$ws.Range("E46").Value = '  -5.51%  '
$ws.Range("D47").Value = '12.01'
$ws.Range("E47").Value = '  -7.36%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.928'
$ws.Range("E48").Value = '  -6.40%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '3.371'
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("E51").Value = '  -2.81%  '
